# Add the "2022-Q1" worksheet, a snapshot of fund holdings of 600699 for
# the new quarter, positioned right before the "总计" (totals) sheet, i.e.
# right after "2021-Q4".

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# NOTE: sheet handles are position-based, so any handle captured before an
# Add/Move/Delete on Worksheets can end up pointing at the wrong sheet
# afterwards. Re-resolve "总计" by name now that the new sheet has been
# inserted ahead of it.
$totals = $wb.Worksheets.Item("总计")

# --- copy the header/index-column formatting from the "2021-Q4" sheet so
# the new sheet matches the existing look (bold + thin border, centered).
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A8").Copy()
$newSheet.Range("A2:A8").PasteSpecial(-4122)

# --- header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- fund holdings data rows (B:G are text even though numeric-looking,
# H is a true number). Build as an array-of-arrays using the unary comma
# operator -- a plain `@(@(...), @(...))` literal flattens into one array
# in this host, which breaks per-row indexing.
$data = @()
$data += , @("562500", "华夏中证机器人ETF", "1.35", "99.22", "2.52", "0.0340", 8)
$data += , @("159770", "天弘中证机器人ETF", "0.87", "99.47", "2.53", "0.0220", 8)
$data += , @("562360", "银华中证机器人ETF", "0.68", "96.94", "2.44", "0.0166", 8)
$data += , @("006346", "安信量化优选股票A", "0.71", "90.62", "0.61", "0.0043", 7)
$data += , @("006347", "安信量化优选股票C", "0.49", "90.62", "0.61", "0.0030", 7)
$data += , @("010999", "兴华瑞丰混合A", "0.06", "29.21", "2.64", "0.0016", 7)
$data += , @("011000", "兴华瑞丰混合C", "0.05", "29.21", "2.64", "0.0013", 7)

# --- index column (A2:A8), numeric 0..6
for ($i = 0; $i -lt $data.Count; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $i
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $textRange = $newSheet.Range("B$row:G$row")
    $textRange.NumberFormat = "@"
    $newSheet.Cells.Item($row, 2).Value = $rec[0]
    $newSheet.Cells.Item($row, 3).Value = $rec[1]
    $newSheet.Cells.Item($row, 4).Value = $rec[2]
    $newSheet.Cells.Item($row, 5).Value = $rec[3]
    $newSheet.Cells.Item($row, 6).Value = $rec[4]
    $newSheet.Cells.Item($row, 7).Value = $rec[5]
    $textRange.ClearFormats()
    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}

# =====================================================================
# Update the "总计" (totals) sheet: insert a new top data row for
# "2022-Q1" and shift the existing quarters' rows down by one, renumbering
# the index column.
# =====================================================================

$totals.Rows.Item(2).Insert()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 7
$totals.Range("D2").Value = 0.08

# renumber the index column for the rows that got shifted down
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5
